$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Target cluster (column D) is "ECs".
# These were original rows 2, 5 and 8; delete bottom-up so row numbers
# of not-yet-deleted rows remain stable.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Refresh derived-specificity / weight columns (G,H,I,J,M,N,O,P,Q,R,S,T)
# for the remaining 6 data rows using the recomputed TPM-based figures.

# Row 2: ECs -> FAPs
$ws.Range("G2").Value = 0.1728506666666667
$ws.Range("H2").Value = 0.518552
$ws.Range("I2").Value = 0.0840503369699626
$ws.Range("J2").Value = 0.0840503369699626
$ws.Range("M2").Value = 0.06103333333333334
$ws.Range("N2").Value = 0.1831
$ws.Range("O2").Value = 0.2094151016766933
$ws.Range("P2").Value = 0.2094151016766933
$ws.Range("Q2").Value = 0.01054965235555556
$ws.Range("R2").Value = 0.0949468712
$ws.Range("S2").Value = 0.01760140986252505
$ws.Range("T2").Value = 0.01760140986252505

# Row 3: ECs -> MuSCs
$ws.Range("G3").Value = 0.1728506666666667
$ws.Range("H3").Value = 0.518552
$ws.Range("I3").Value = 0.0840503369699626
$ws.Range("J3").Value = 0.0840503369699626
$ws.Range("M3").Value = 0.2304133333333333
$ws.Range("N3").Value = 0.69124
$ws.Range("O3").Value = 0.7905848983233067
$ws.Range("P3").Value = 0.7905848983233067
$ws.Range("Q3").Value = 0.03982709827555556
$ws.Range("R3").Value = 0.35844388448
$ws.Range("S3").Value = 0.06644892710743755
$ws.Range("T3").Value = 0.06644892710743755

# Row 4: FAPs -> FAPs
$ws.Range("G4").Value = 1.367717666666667
$ws.Range("H4").Value = 4.103153
$ws.Range("I4").Value = 0.6650661694281633
$ws.Range("J4").Value = 0.6650661694281633
$ws.Range("M4").Value = 0.06103333333333334
$ws.Range("N4").Value = 0.1831
$ws.Range("O4").Value = 0.2094151016766933
$ws.Range("P4").Value = 0.2094151016766933
$ws.Range("Q4").Value = 0.08347636825555556
$ws.Range("R4").Value = 0.7512873143000001
$ws.Range("S4").Value = 0.1392748994925277
$ws.Range("T4").Value = 0.1392748994925277

# Row 5: FAPs -> MuSCs
$ws.Range("G5").Value = 1.367717666666667
$ws.Range("H5").Value = 4.103153
$ws.Range("I5").Value = 0.6650661694281633
$ws.Range("J5").Value = 0.6650661694281633
$ws.Range("M5").Value = 0.2304133333333333
$ws.Range("N5").Value = 0.69124
$ws.Range("O5").Value = 0.7905848983233067
$ws.Range("P5").Value = 0.7905848983233067
$ws.Range("Q5").Value = 0.3151403866355555
$ws.Range("R5").Value = 2.83626347972
$ws.Range("S5").Value = 0.5257912699356355
$ws.Range("T5").Value = 0.5257912699356355

# Row 6: MuSCs -> FAPs
$ws.Range("G6").Value = 0.5159453333333334
$ws.Range("H6").Value = 1.547836
$ws.Range("I6").Value = 0.2508834936018741
$ws.Range("J6").Value = 0.2508834936018741
$ws.Range("M6").Value = 0.06103333333333334
$ws.Range("N6").Value = 0.1831
$ws.Range("O6").Value = 0.2094151016766933
$ws.Range("P6").Value = 0.2094151016766933
$ws.Range("Q6").Value = 0.03148986351111111
$ws.Range("R6").Value = 0.2834087716
$ws.Range("S6").Value = 0.05253879232164049
$ws.Range("T6").Value = 0.05253879232164049

# Row 7: MuSCs -> MuSCs
$ws.Range("G7").Value = 0.5159453333333334
$ws.Range("H7").Value = 1.547836
$ws.Range("I7").Value = 0.2508834936018741
$ws.Range("J7").Value = 0.2508834936018741
$ws.Range("M7").Value = 0.2304133333333333
$ws.Range("N7").Value = 0.69124
$ws.Range("O7").Value = 0.7905848983233067
$ws.Range("P7").Value = 0.7905848983233067
$ws.Range("Q7").Value = 0.1188806840711111
$ws.Range("R7").Value = 1.06992615664
$ws.Range("S7").Value = 0.1983447012802336
$ws.Range("T7").Value = 0.1983447012802336
